$d = $word.ActiveDocument

# Helper: find the (1-based) index of the paragraph whose text equals
# $text. Always re-reads from the live $d.Paragraphs collection (never
# caches a Paragraph object across a structural edit), because once
# paragraphs are inserted/removed, previously-fetched Paragraph
# references stop tracking their live position in this host.
function Find-ParaIndex([string]$text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Add a new "Version 4.0" section (a "Layers added" bullet point,
#    then a blank line) right before the existing "Version 3.0"
#    section.
# ------------------------------------------------------------------
$idx3 = Find-ParaIndex "Version 3.0"

$p3 = $d.Paragraphs($idx3)
$p3.Range.InsertParagraphBefore()
$p3.Range.InsertParagraphBefore()
$p3.Range.InsertParagraphBefore()

# The three new blank paragraphs now occupy $idx3, $idx3+1, $idx3+2;
# "Version 3.0" itself shifted to $idx3+3.
$d.Paragraphs($idx3).Range.Text = "Version 4.0"
$d.Paragraphs($idx3 + 1).Range.Text = "Layers added"
# $d.Paragraphs($idx3 + 2) is left blank on purpose, separating the
# new section from "Version 3.0" below it.

# ------------------------------------------------------------------
# 2) Move the (hidden) "_GoBack" bookmark. It currently sits in the
#    middle of the "Glitch removed: Background images..." paragraph,
#    splitting it into two runs. Move it down so that it instead
#    occupies its own paragraph, right after "Undo function now
#    added" (taking the place of the blank paragraph that used to
#    follow it there).
# ------------------------------------------------------------------
$bookmark = $d.Bookmarks.Item("_GoBack")
$bookmark.Delete()

$idxUndo = Find-ParaIndex "Undo function now added"
$blankAfterUndo = $d.Paragraphs($idxUndo + 1)
$blankAfterUndo.Range.Bookmarks.Add("_GoBack")
